$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 1005
$ws.Range("I106").Value = 1005
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 1005
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -374
$ws.Range("N106").ClearContents()

$ws.Range("H137").Value = 41676460
$ws.Range("I137").Value = 925
$ws.Range("J137").Value = 50011564
$ws.Range("K137").Value = 2775
$ws.Range("L137").Value = 150034692
$ws.Range("M137").Value = -225
$ws.Range("N137").Value = -150039792

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1962153.6
$ws.Range("I2").Value = 1825
$ws.Range("J2").Value = 4202529
$ws.Range("K2").Value = 1825
$ws.Range("L2").Value = 4202529
$ws.Range("M2").Value = -1712
$ws.Range("N2").Value = -4202755

$ws.Range("H32").Value = 5365.721
$ws.Range("I32").Value = 6529.595
$ws.Range("J32").Value = 2792.9473
$ws.Range("K32").Value = 6529.595
$ws.Range("L32").Value = 2792.9473
$ws.Range("M32").Value = -6242.595
$ws.Range("N32").Value = -3366.9473

$ws.Range("H61").Value = 2585.7585
$ws.Range("I61").Value = 1782.0526
$ws.Range("J61").Value = 4112.8
$ws.Range("K61").Value = 1782.0526
$ws.Range("L61").Value = 4112.8
$ws.Range("M61").Value = -1570.0526
$ws.Range("N61").Value = -4536.8

$ws.Range("H74").Value = 855.0714
$ws.Range("I74").Value = 500.125
$ws.Range("J74").Value = 1328.3334
$ws.Range("K74").Value = 500.125
$ws.Range("L74").Value = 1328.3334
$ws.Range("M74").Value = 373.875
$ws.Range("N74").Value = -3076.3334

$ws.Range("H77").Value = 855.0714
$ws.Range("I77").Value = 500.125
$ws.Range("J77").Value = 1328.3334
$ws.Range("K77").Value = 2500.625
$ws.Range("L77").Value = 6641.666999999999
$ws.Range("M77").Value = 1867.375
$ws.Range("N77").Value = -15377.667

$ws.Range("H112").Value = 31230.834
$ws.Range("J112").Value = 31230.834
$ws.Range("L112").Value = 31230.834
$ws.Range("N112").Value = -34184.834

$ws.Range("H116").Value = 1962153.6
$ws.Range("I116").Value = 1825
$ws.Range("J116").Value = 4202529
$ws.Range("K116").Value = 1825
$ws.Range("L116").Value = 4202529
$ws.Range("M116").Value = 469
$ws.Range("N116").Value = -4207117

$ws.Range("H132").Value = 2483.6177
$ws.Range("I132").Value = 2371.9614
$ws.Range("J132").Value = 2846.5
$ws.Range("K132").Value = 7115.8842
$ws.Range("L132").Value = 8539.5
$ws.Range("M132").Value = -4585.8842
$ws.Range("N132").Value = -13599.5

$ws.Range("H136").Value = 2585.7585
$ws.Range("I136").Value = 1782.0526
$ws.Range("J136").Value = 4112.8
$ws.Range("K136").Value = 5346.1578
$ws.Range("L136").Value = 12338.4
$ws.Range("M136").Value = -2796.1578
$ws.Range("N136").Value = -17438.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1962153.6
$ws.Range("I3").Value = 1825
$ws.Range("J3").Value = 4202529
$ws.Range("K3").Value = 1825
$ws.Range("L3").Value = 4202529
$ws.Range("M3").Value = -1711
$ws.Range("N3").Value = -4202757

$ws.Range("H86").Value = 1918.2333
$ws.Range("I86").Value = 1689.96
$ws.Range("K86").Value = 1689.96
$ws.Range("M86").Value = -566.96

$ws.Range("H89").Value = 1918.2333
$ws.Range("I89").Value = 1689.96
$ws.Range("K89").Value = 8449.799999999999
$ws.Range("M89").Value = -2833.799999999999

$ws.Range("H134").Value = 9938.5
$ws.Range("I134").Value = 10660.182
$ws.Range("K134").Value = 31980.546
$ws.Range("M134").Value = -29445.546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1722.4584
$ws.Range("I31").Value = 1063.762
$ws.Range("J31").Value = 6333.3335
$ws.Range("K31").Value = 1063.762
$ws.Range("L31").Value = 6333.3335
$ws.Range("M31").Value = -768.7619999999999
$ws.Range("N31").Value = -6923.3335

$ws.Range("H34").Value = 1722.4584
$ws.Range("I34").Value = 1063.762
$ws.Range("J34").Value = 6333.3335
$ws.Range("K34").Value = 1063.762
$ws.Range("L34").Value = 6333.3335
$ws.Range("M34").Value = -861.7619999999999
$ws.Range("N34").Value = -6737.3335

$ws.Range("H58").Value = 4381.5386
$ws.Range("I58").Value = 3383.1428
$ws.Range("J58").Value = 4749.3687
$ws.Range("K58").Value = 3383.1428
$ws.Range("L58").Value = 4749.3687
$ws.Range("M58").Value = -3180.1428
$ws.Range("N58").Value = -5155.3687

$ws.Range("H81").Value = 50328
$ws.Range("J81").Value = 50328
$ws.Range("L81").Value = 50328
$ws.Range("N81").Value = -52324

$ws.Range("H84").Value = 50328
$ws.Range("J84").Value = 50328
$ws.Range("L84").Value = 150984
$ws.Range("N84").Value = -160968

$ws.Range("H132").Value = 4461.8237
$ws.Range("I132").Value = 3604.1538
$ws.Range("J132").Value = 7249.25
$ws.Range("K132").Value = 10812.4614
$ws.Range("L132").Value = 21747.75
$ws.Range("M132").Value = -8282.4614
$ws.Range("N132").Value = -26807.75

$ws.Range("H134").Value = 1811.3704
$ws.Range("I134").Value = 1876.28
$ws.Range("K134").Value = 5628.84
$ws.Range("M134").Value = -3093.84

$ws.Range("H136").Value = 4381.5386
$ws.Range("I136").Value = 3383.1428
$ws.Range("J136").Value = 4749.3687
$ws.Range("K136").Value = 10149.4284
$ws.Range("L136").Value = 14248.1061
$ws.Range("M136").Value = -7599.428400000001
$ws.Range("N136").Value = -19348.1061

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 38462268
$ws.Range("J121").Value = 41667430
$ws.Range("L121").Value = 125002290
$ws.Range("N121").Value = -125004910

$ws.Range("H122").Value = 18519896
$ws.Range("I122").Value = 41667164
$ws.Range("J122").Value = 2082.9
$ws.Range("K122").Value = 375004476
$ws.Range("L122").Value = 18746.1
$ws.Range("M122").Value = -375002026
$ws.Range("N122").Value = -23646.1

$ws.Range("H131").Value = 2520.1792
$ws.Range("I131").Value = 3243.2
$ws.Range("J131").Value = 2212.5107
$ws.Range("K131").Value = 9729.599999999999
$ws.Range("L131").Value = 6637.532099999999
$ws.Range("M131").Value = -4689.599999999999
$ws.Range("N131").Value = -16717.5321

$ws.Range("H133").Value = 6653.3184
$ws.Range("J133").Value = 7057
$ws.Range("L133").Value = 21171
$ws.Range("N133").Value = -31291

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 21250
$ws.Range("J110").Value = 21250
$ws.Range("L110").Value = 21250
$ws.Range("N110").Value = -29430

$ws.Range("H136").Value = 1834.375
$ws.Range("I136").Value = 1492.2
$ws.Range("J136").Value = 2078.7856
$ws.Range("K136").Value = 4476.6
$ws.Range("L136").Value = 6236.3568
$ws.Range("M136").Value = -1926.6
$ws.Range("N136").Value = -11336.3568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 34000
$ws.Range("J75").Value = 34000
$ws.Range("L75").Value = 34000
$ws.Range("N75").Value = -35872

$ws.Range("H78").Value = 34000
$ws.Range("J78").Value = 34000
$ws.Range("L78").Value = 102000
$ws.Range("N78").Value = -111360

$ws.Range("H122").Value = 961.9375
$ws.Range("I122").Value = 976.61536
$ws.Range("J122").Value = 898.3333
$ws.Range("K122").Value = 2929.84608
$ws.Range("L122").Value = 2694.9999
$ws.Range("M122").Value = -479.8460800000003
$ws.Range("N122").Value = -7594.9999
